$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 11.00333457712453
$ws.Range("C2").Value = 3.568645991273135
$ws.Range("D2").Value = 8.468525394246994
$ws.Range("F2").Value = 41.64567160887763
$ws.Range("G2").Value = 48.56230286015465
$ws.Range("H2").Value = 19.13241850530223
$ws.Range("J2").Value = 11.42487781406987
$ws.Range("K2").Value = 10.47579132728361
$ws.Range("M2").Value = 16.41127631392882
$ws.Range("B3").Value = 10.82338134992216
$ws.Range("C3").Value = 3.384346221411306
$ws.Range("D3").Value = 8.455276172373877
$ws.Range("F3").Value = 41.59500022875628
$ws.Range("G3").Value = 48.42690861508594
$ws.Range("H3").Value = 19.15162029634745
$ws.Range("J3").Value = 11.43965351529483
$ws.Range("K3").Value = 10.37135141627131
$ws.Range("M3").Value = 16.37467763700988
$ws.Range("B4").Value = 10.71448951913113
$ws.Range("C4").Value = 3.265096089946604
$ws.Range("D4").Value = 8.448714176122843
$ws.Range("F4").Value = 41.57214417567096
$ws.Range("G4").Value = 48.35447850577844
$ws.Range("H4").Value = 19.16625869243896
$ws.Range("J4").Value = 11.45027588544444
$ws.Range("K4").Value = 10.30939957790457
$ws.Range("M4").Value = 16.35549117076285
$ws.Range("B5").Value = 10.67058129116512
$ws.Range("C5").Value = 3.214992042616315
$ws.Range("D5").Value = 8.446438277645102
$ws.Range("F5").Value = 41.56491072644432
$ws.Range("G5").Value = 48.32766845995567
$ws.Range("H5").Value = 19.17293963779389
$ws.Range("J5").Value = 11.45499441938526
$ws.Range("K5").Value = 10.28472928169911
$ws.Range("M5").Value = 16.34850488110603
$ws.Range("B6").Value = 10.66332036701215
$ws.Range("C6").Value = 3.206581939707876
$ws.Range("D6").Value = 8.446084485521203
$ws.Range("F6").Value = 41.56383538043853
$ws.Range("G6").Value = 48.3233804651949
$ws.Range("H6").Value = 19.17409220927617
$ws.Range("J6").Value = 11.45580147416297
$ws.Range("K6").Value = 10.2806683789844
$ws.Range("M6").Value = 16.34739525331963
$ws.Range("B7").Value = 10.71389538640994
$ws.Range("C7").Value = 3.264426444374703
$ws.Range("D7").Value = 8.448681867150936
$ws.Range("F7").Value = 41.57203819396887
$ws.Range("G7").Value = 48.35410596341775
$ws.Range("H7").Value = 19.16634589715336
$ws.Range("J7").Value = 11.45033794278217
$ws.Range("K7").Value = 10.30906449803826
$ws.Range("M7").Value = 16.35539357335239
$ws.Range("B8").Value = 10.94099288040755
$ws.Range("C8").Value = 3.506375369550395
$ws.Range("D8").Value = 8.4636322129495
$ws.Range("F8").Value = 41.62648959658711
$ws.Range("G8").Value = 48.51340722962434
$ws.Range("H8").Value = 19.13844778433698
$ws.Range("J8").Value = 11.42965081840075
$ws.Range("K8").Value = 10.43934666421322
$ws.Range("M8").Value = 16.39797918955622
$ws.Range("B9").Value = 11.39622610906249
$ws.Range("C9").Value = 3.931752569814363
$ws.Range("D9").Value = 8.505305753161295
$ws.Range("F9").Value = 41.79852546690216
$ws.Range("G9").Value = 48.90994909307432
$ws.Range("H9").Value = 19.10636263823905
$ws.Range("J9").Value = 11.40137849069758
$ws.Range("K9").Value = 10.71077219207479
$ws.Range("M9").Value = 16.50725176707424
$ws.Range("B10").Value = 11.73300776210445
$ws.Range("C10").Value = 4.213584215782403
$ws.Range("D10").Value = 8.543271061279539
$ws.Range("F10").Value = 41.96427694745602
$ws.Range("G10").Value = 49.25138247747304
$ws.Range("H10").Value = 19.09660571773447
$ws.Range("J10").Value = 11.38809577993663
$ws.Range("K10").Value = 10.91809463200976
$ws.Range("M10").Value = 16.60278776070343
$ws.Range("B11").Value = 11.88596120307514
$ws.Range("C11").Value = 4.335028105926777
$ws.Range("D11").Value = 8.562093047918438
$ws.Range("F11").Value = 42.04810598802722
$ws.Range("G11").Value = 49.41727128808733
$ws.Range("H11").Value = 19.09516924593602
$ws.Range("J11").Value = 11.38367690732641
$ws.Range("K11").Value = 11.01373257015708
$ws.Range("M11").Value = 16.6494512455423
$ws.Range("B12").Value = 11.94378385244831
$ws.Range("C12").Value = 4.380037330348477
$ws.Range("D12").Value = 8.569439301269949
$ws.Range("F12").Value = 42.08104889676069
$ws.Range("G12").Value = 49.48157783982714
$ws.Range("H12").Value = 19.09505680967186
$ws.Range("J12").Value = 11.382236732615
$ws.Range("K12").Value = 11.05010554420739
$ws.Range("M12").Value = 16.66757178797112
$ws.Range("B13").Value = 11.93133599992544
$ws.Range("C13").Value = 4.370387411059539
$ws.Range("D13").Value = 8.567847495365884
$ws.Range("F13").Value = 42.0739009673375
$ws.Range("G13").Value = 49.4676626240709
$ws.Range("H13").Value = 19.09506183706891
$ws.Range("J13").Value = 11.3825365357554
$ws.Range("K13").Value = 11.04226552720354
$ws.Range("M13").Value = 16.66364936208776
$ws.Range("B14").Value = 11.8907206049059
$ws.Range("C14").Value = 4.338750691959644
$ws.Range("D14").Value = 8.562693069595745
$ws.Range("F14").Value = 42.0507922789365
$ws.Range("G14").Value = 49.42253219480502
$ws.Range("H14").Value = 19.09515134888482
$ws.Range("J14").Value = 11.383553752425
$ws.Range("K14").Value = 11.01672204328654
$ws.Range("M14").Value = 16.65093307578337
$ws.Range("B15").Value = 11.8658280213786
$ws.Range("C15").Value = 4.319244640137879
$ws.Range("D15").Value = 8.559564193797028
$ws.Range("F15").Value = 42.03679322395547
$ws.Range("G15").Value = 49.39508130705521
$ws.Range("H15").Value = 19.09526236693872
$ws.Range("J15").Value = 11.38420718088355
$ws.Range("K15").Value = 11.00109536555695
$ws.Range("M15").Value = 16.64320227383145
$ws.Range("B16").Value = 11.72300154638999
$ws.Range("C16").Value = 4.205510893005307
$ws.Range("D16").Value = 8.542071855143629
$ws.Range("F16").Value = 41.95896701421849
$ws.Range("G16").Value = 49.2407512738101
$ws.Range("H16").Value = 19.09675999282513
$ws.Range("J16").Value = 11.38841720829786
$ws.Range("K16").Value = 10.91186826455129
$ws.Range("M16").Value = 16.59980176724298
$ws.Range("B17").Value = 11.63527352138451
$ws.Range("C17").Value = 4.134001620014343
$ws.Range("D17").Value = 8.531735230468243
$ws.Range("F17").Value = 41.91337278104957
$ws.Range("G17").Value = 49.14876000997302
$ws.Range("H17").Value = 19.09844756416298
$ws.Range("J17").Value = 11.39141554626942
$ws.Range("K17").Value = 10.85744591226331
$ws.Range("M17").Value = 16.57399016436453
$ws.Range("B18").Value = 11.58479466816822
$ws.Range("C18").Value = 4.092235092621005
$ws.Range("D18").Value = 8.525936174686276
$ws.Range("F18").Value = 41.88794234199559
$ws.Range("G18").Value = 49.09684594518998
$ws.Range("H18").Value = 19.09970079349473
$ws.Range("J18").Value = 11.39329294423067
$ws.Range("K18").Value = 10.82627043161821
$ws.Range("M18").Value = 16.55944623323466
$ws.Range("B19").Value = 11.56770187196122
$ws.Range("C19").Value = 4.077984609352647
$ws.Range("D19").Value = 8.523997967942581
$ws.Range("F19").Value = 41.87946880900112
$ws.Range("G19").Value = 49.07944087499528
$ws.Range("H19").Value = 19.10017365134279
$ws.Range("J19").Value = 11.39395485634738
$ws.Range("K19").Value = 10.81573775946322
$ws.Range("M19").Value = 16.55457412718568
$ws.Range("B20").Value = 11.64461484402002
$ws.Range("C20").Value = 4.141679809983517
$ws.Range("D20").Value = 8.532820471016258
$ws.Range("F20").Value = 41.91814426381207
$ws.Range("G20").Value = 49.15844968583448
$ws.Range("H20").Value = 19.0982386736025
$ws.Range("J20").Value = 11.39108055238682
$ws.Range("K20").Value = 10.86322638408472
$ws.Range("M20").Value = 16.57670664965778
$ws.Range("B21").Value = 11.90265345514547
$ws.Range("C21").Value = 4.348069773993397
$ws.Range("D21").Value = 8.564201147479647
$ws.Range("F21").Value = 42.05754744219038
$ws.Range("G21").Value = 49.4357479760248
$ws.Range("H21").Value = 19.09511334778594
$ws.Range("J21").Value = 11.38324864613441
$ws.Range("K21").Value = 11.0242207915081
$ws.Range("M21").Value = 16.65465602936931
$ws.Range("B22").Value = 12.07070155895823
$ws.Range("C22").Value = 4.477249744721286
$ws.Range("D22").Value = 8.585983251079385
$ws.Range("F22").Value = 42.15563506822831
$ws.Range("G22").Value = 49.62563549625077
$ws.Range("H22").Value = 19.09558582282476
$ws.Range("J22").Value = 11.37948887371846
$ws.Range("K22").Value = 11.13033937295328
$ws.Range("M22").Value = 16.70821860539717
$ws.Range("B23").Value = 11.98108524637988
$ws.Range("C23").Value = 4.408827909494788
$ws.Range("D23").Value = 8.574242739397675
$ws.Range("F23").Value = 42.10264997885585
$ws.Range("G23").Value = 49.52350797224977
$ws.Range("H23").Value = 19.09510362446879
$ws.Range("J23").Value = 11.38137131717275
$ws.Range("K23").Value = 11.07363067092213
$ws.Range("M23").Value = 16.6793953457372
$ws.Range("B24").Value = 11.6403917607779
$ws.Range("C24").Value = 4.138210539685884
$ws.Range("D24").Value = 8.532329386128273
$ws.Range("F24").Value = 41.91598463857139
$ws.Range("G24").Value = 49.15406595022822
$ws.Range("H24").Value = 19.09833223135515
$ws.Range("J24").Value = 11.39123152460367
$ws.Range("K24").Value = 10.86061267928489
$ws.Range("M24").Value = 16.57547760538998
$ws.Range("B25").Value = 11.27240142346442
$ws.Range("C25").Value = 3.822023165791578
$ws.Range("D25").Value = 8.492727149726868
$ws.Range("F25").Value = 41.74503691621727
$ws.Range("G25").Value = 48.79377822164401
$ws.Range("H25").Value = 19.11261741371706
$ws.Range("J25").Value = 11.40771107697538
$ws.Range("K25").Value = 10.63582258173244
$ws.Range("M25").Value = 16.47497648994393
